$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.219.52"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.631.67"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.95"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.84"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.630.70"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.111.51"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.142.34"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.75"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.632.25"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.08"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.55"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.82"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.767.38"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "494.31"
$ws.Range("E32").Value = "  -5.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.26"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.24"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.34"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -7.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.86"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.04"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.55"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.64"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.543"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.599"
$ws.Range("E51").Value = "  -0.53%  "
